{"js": "const tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// (row, col) -> new text, using 0-based indices into the full 20-row x 5-col table\n// (rows 0,4,8,12,16 hold the division problems; the rows between them are blank spacers).\nconst updates = [\n  { row: 0, col: 0, before: \"57\u00f77=8, 1\", after: \"31\u00f79=3, 4\" },\n  { row: 0, col: 1, before: \"73\u00f74=18, 1\", after: \"22\u00f72=11, 0\" },\n  { row: 0, col: 2, before: \"28\u00f73=9, 1\", after: \"18\u00f73=6, 0\" },\n  { row: 0, col: 3, before: \"50\u00f79=5, 5\", after: \"27\u00f72=13, 1\" },\n  { row: 0, col: 4, before: \"96\u00f72=48, 0\", after: \"88\u00f73=29, 1\" },\n  { row: 4, col: 0, before: \"77\u00f79=8, 5\", after: \"76\u00f78=9, 4\" },\n  { row: 4, col: 1, before: \"60\u00f75=12, 0\", after: \"41\u00f79=4, 5\" },\n  { row: 4, col: 2, before: \"48\u00f75=9, 3\", after: \"68\u00f77=9, 5\" },\n  { row: 4, col: 3, before: \"60\u00f72=30, 0\", after: \"38\u00f79=4, 2\" },\n  { row: 4, col: 4, before: \"41\u00f76=6, 5\", after: \"16\u00f73=5, 1\" },\n  { row: 8, col: 0, before: \"12\u00f75=2, 2\", after: \"69\u00f78=8, 5\" },\n  { row: 8, col: 1, before: \"63\u00f72=31, 1\", after: \"58\u00f79=6, 4\" },\n  { row: 8, col: 2, before: \"70\u00f78=8, 6\", after: \"19\u00f74=4, 3\" },\n  { row: 8, col: 3, before: \"36\u00f79=4, 0\", after: \"38\u00f74=9, 2\" },\n  { row: 8, col: 4, before: \"22\u00f72=11, 0\", after: \"93\u00f75=18, 3\" },\n  { row: 12, col: 0, before: \"45\u00f75=9, 0\", after: \"30\u00f75=6, 0\" },\n  { row: 12, col: 1, before: \"13\u00f76=2, 1\", after: \"16\u00f76=2, 4\" },\n  { row: 12, col: 2, before: \"97\u00f79=10, 7\", after: \"47\u00f75=9, 2\" },\n  { row: 12, col: 3, before: \"27\u00f73=9, 0\", after: \"81\u00f78=10, 1\" },\n  { row: 12, col: 4, before: \"45\u00f72=22, 1\", after: \"42\u00f74=10, 2\" },\n  { row: 16, col: 0, before: \"86\u00f72=43, 0\", after: \"84\u00f78=10, 4\" },\n  { row: 16, col: 1, before: \"72\u00f78=9, 0\", after: \"34\u00f78=4, 2\" },\n  { row: 16, col: 2, before: \"11\u00f72=5, 1\", after: \"17\u00f78=2, 1\" },\n  { row: 16, col: 3, before: \"34\u00f72=17, 0\", after: \"55\u00f77=7, 6\" },\n  { row: 16, col: 4, before: \"12\u00f79=1, 3\", after: \"29\u00f77=4, 1\" },\n];\n\nconst cells = updates.map(u => table.getCell(u.row, u.col));\ncells.forEach(c => c.load(\"value\"));\nawait context.sync();\n\nupdates.forEach((u, i) => {\n  const cell = cells[i];\n  if (cell.value === u.before) {\n    cell.value = u.after;\n  } else if (cell.value !== u.after) {\n    // Unexpected existing content - only overwrite if it still matches the\n    // pre-edit text; otherwise leave alone to avoid clobbering other data.\n    throw new Error(\n      `Unexpected text in cell (${u.row}, ${u.col}): ${JSON.stringify(cell.value)}`\n    );\n  }\n});\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# (row, col) use 1-based indices matching Word's Table.Cell(row, col); the full\n# table has 20 rows x 5 columns - rows 1,5,9,13,17 hold the division problems,\n# the rows in between are blank spacer rows.\n$updates = @(\n    @{ Row = 1; Col = 1; Before = \"57\u00f77=8, 1\"; After = \"31\u00f79=3, 4\" }\n    @{ Row = 1; Col = 2; Before = \"73\u00f74=18, 1\"; After = \"22\u00f72=11, 0\" }\n    @{ Row = 1; Col = 3; Before = \"28\u00f73=9, 1\"; After = \"18\u00f73=6, 0\" }\n    @{ Row = 1; Col = 4; Before = \"50\u00f79=5, 5\"; After = \"27\u00f72=13, 1\" }\n    @{ Row = 1; Col = 5; Before = \"96\u00f72=48, 0\"; After = \"88\u00f73=29, 1\" }\n    @{ Row = 5; Col = 1; Before = \"77\u00f79=8, 5\"; After = \"76\u00f78=9, 4\" }\n    @{ Row = 5; Col = 2; Before = \"60\u00f75=12, 0\"; After = \"41\u00f79=4, 5\" }\n    @{ Row = 5; Col = 3; Before = \"48\u00f75=9, 3\"; After = \"68\u00f77=9, 5\" }\n    @{ Row = 5; Col = 4; Before = \"60\u00f72=30, 0\"; After = \"38\u00f79=4, 2\" }\n    @{ Row = 5; Col = 5; Before = \"41\u00f76=6, 5\"; After = \"16\u00f73=5, 1\" }\n    @{ Row = 9; Col = 1; Before = \"12\u00f75=2, 2\"; After = \"69\u00f78=8, 5\" }\n    @{ Row = 9; Col = 2; Before = \"63\u00f72=31, 1\"; After = \"58\u00f79=6, 4\" }\n    @{ Row = 9; Col = 3; Before = \"70\u00f78=8, 6\"; After = \"19\u00f74=4, 3\" }\n    @{ Row = 9; Col = 4; Before = \"36\u00f79=4, 0\"; After = \"38\u00f74=9, 2\" }\n    @{ Row = 9; Col = 5; Before = \"22\u00f72=11, 0\"; After = \"93\u00f75=18, 3\" }\n    @{ Row = 13; Col = 1; Before = \"45\u00f75=9, 0\"; After = \"30\u00f75=6, 0\" }\n    @{ Row = 13; Col = 2; Before = \"13\u00f76=2, 1\"; After = \"16\u00f76=2, 4\" }\n    @{ Row = 13; Col = 3; Before = \"97\u00f79=10, 7\"; After = \"47\u00f75=9, 2\" }\n    @{ Row = 13; Col = 4; Before = \"27\u00f73=9, 0\"; After = \"81\u00f78=10, 1\" }\n    @{ Row = 13; Col = 5; Before = \"45\u00f72=22, 1\"; After = \"42\u00f74=10, 2\" }\n    @{ Row = 17; Col = 1; Before = \"86\u00f72=43, 0\"; After = \"84\u00f78=10, 4\" }\n    @{ Row = 17; Col = 2; Before = \"72\u00f78=9, 0\"; After = \"34\u00f78=4, 2\" }\n    @{ Row = 17; Col = 3; Before = \"11\u00f72=5, 1\"; After = \"17\u00f78=2, 1\" }\n    @{ Row = 17; Col = 4; Before = \"34\u00f72=17, 0\"; After = \"55\u00f77=7, 6\" }\n    @{ Row = 17; Col = 5; Before = \"12\u00f79=1, 3\"; After = \"29\u00f77=4, 1\" }\n)\n\nforeach ($u in $updates) {\n    $cell = $table.Cell($u.Row, $u.Col)\n    $range = $cell.Range\n    # Cell range text includes trailing cell-mark characters; trim them for comparison.\n    $current = $range.Text.TrimEnd([char]7, [char]13)\n    if ($current -eq $u.Before) {\n        $range.Text = $u.After\n    } elseif ($current -ne $u.After) {\n        throw \"Unexpected text in cell ($($u.Row), $($u.Col)): $current\"\n    }\n}\n"}
